# Applies updated crypto price/volume figures to Sheet1 (cryptos.xlsx)
# Commit: "Updated cryptos list on Wed May  8 11:57:22 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "62.292.32"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.993.69"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.06"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.32"
$ws.Range("E6").Value = "  -6.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "2.994.08"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.45"
$ws.Range("E14").Value = "  -5.68%  "
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "3.488.58"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "62.275.82"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "2.989.65"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.55"
$ws.Range("E20").Value = "  -4.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.80"
$ws.Range("E21").Value = "  -3.89%  "
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.22"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("E25").Value = "  -7.49%  "
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  -6.11%  "
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.75"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  -6.09%  "
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.03"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.91"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  -11.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "392.47"
$ws.Range("E43").Value = "  -10.17%  "
$ws.Range("E44").Value = "  -5.26%  "
$ws.Range("D45").Value = "2.757.59"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.85"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.34"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.66"
$ws.Range("E51").Value = "  -7.75%  "
